# Apply updated NATMI Nid1-Itgb3 ligand-receptor results ("Natmi following Dr Hou advice").
# Sending/target cluster order is now ECs, FAPs, M2, sCs (4x4 = 16 data rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Nid1"
$ws.Cells.Item(2,3).Value = "Itgb3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 55.41713066666667
$ws.Cells.Item(2,8).Value = 166.251392
$ws.Cells.Item(2,9).Value = 0.08138603925734667
$ws.Cells.Item(2,10).Value = 0.08138603925734668
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 3.778439
$ws.Cells.Item(2,14).Value = 11.335317
$ws.Cells.Item(2,15).Value = 0.4252971528324392
$ws.Cells.Item(2,16).Value = 0.4252971528324392
$ws.Cells.Item(2,17).Value = 209.3902477790294
$ws.Cells.Item(2,18).Value = 1884.512230011264
$ws.Cells.Item(2,19).Value = 0.03461325077645867
$ws.Cells.Item(2,20).Value = 0.03461325077645867
# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Nid1"
$ws.Cells.Item(3,3).Value = "Itgb3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 55.41713066666667
$ws.Cells.Item(3,8).Value = 166.251392
$ws.Cells.Item(3,9).Value = 0.08138603925734667
$ws.Cells.Item(3,10).Value = 0.08138603925734668
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 4.333403333333333
$ws.Cells.Item(3,14).Value = 13.00021
$ws.Cells.Item(3,15).Value = 0.4877633593505858
$ws.Cells.Item(3,16).Value = 0.4877633593505858
$ws.Cells.Item(3,17).Value = 240.1447787547022
$ws.Cells.Item(3,18).Value = 2161.30300879232
$ws.Cells.Item(3,19).Value = 0.03969712791240206
$ws.Cells.Item(3,20).Value = 0.03969712791240208
# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Nid1"
$ws.Cells.Item(4,3).Value = "Itgb3"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 55.41713066666667
$ws.Cells.Item(4,8).Value = 166.251392
$ws.Cells.Item(4,9).Value = 0.08138603925734667
$ws.Cells.Item(4,10).Value = 0.08138603925734668
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.2909853333333334
$ws.Cells.Item(4,14).Value = 0.8729560000000001
$ws.Cells.Item(4,15).Value = 0.03275300561492853
$ws.Cells.Item(4,16).Value = 0.03275300561492853
$ws.Cells.Item(4,17).Value = 16.12557223941689
$ws.Cells.Item(4,18).Value = 145.130150154752
$ws.Cells.Item(4,19).Value = 0.002665637400772669
$ws.Cells.Item(4,20).Value = 0.00266563740077267
# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Nid1"
$ws.Cells.Item(5,3).Value = "Itgb3"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 55.41713066666667
$ws.Cells.Item(5,8).Value = 166.251392
$ws.Cells.Item(5,9).Value = 0.08138603925734667
$ws.Cells.Item(5,10).Value = 0.08138603925734668
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.4814053333333333
$ws.Cells.Item(5,14).Value = 1.444216
$ws.Cells.Item(5,15).Value = 0.0541864822020464
$ws.Cells.Item(5,16).Value = 0.05418648220204641
$ws.Cells.Item(5,17).Value = 26.67810226096356
$ws.Cells.Item(5,18).Value = 240.102920348672
$ws.Cells.Item(5,19).Value = 0.004410023167713265
$ws.Cells.Item(5,20).Value = 0.004410023167713266
# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Nid1"
$ws.Cells.Item(6,3).Value = "Itgb3"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 597.374756
$ws.Cells.Item(6,8).Value = 1792.124268
$ws.Cells.Item(6,9).Value = 0.8773093221949784
$ws.Cells.Item(6,10).Value = 0.8773093221949785
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.778439
$ws.Cells.Item(6,14).Value = 11.335317
$ws.Cells.Item(6,15).Value = 0.4252971528324392
$ws.Cells.Item(6,16).Value = 0.4252971528324392
$ws.Cells.Item(6,17).Value = 2257.144075685884
$ws.Cells.Item(6,18).Value = 20314.29668117296
$ws.Cells.Item(6,19).Value = 0.3731171568828814
$ws.Cells.Item(6,20).Value = 0.3731171568828814
# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Nid1"
$ws.Cells.Item(7,3).Value = "Itgb3"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 597.374756
$ws.Cells.Item(7,8).Value = 1792.124268
$ws.Cells.Item(7,9).Value = 0.8773093221949784
$ws.Cells.Item(7,10).Value = 0.8773093221949785
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 4.333403333333333
$ws.Cells.Item(7,14).Value = 13.00021
$ws.Cells.Item(7,15).Value = 0.4877633593505858
$ws.Cells.Item(7,16).Value = 0.4877633593505858
$ws.Cells.Item(7,17).Value = 2588.665758899587
$ws.Cells.Item(7,18).Value = 23297.99183009628
$ws.Cells.Item(7,19).Value = 0.4279193421834081
$ws.Cells.Item(7,20).Value = 0.4279193421834082
# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Nid1"
$ws.Cells.Item(8,3).Value = "Itgb3"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 597.374756
$ws.Cells.Item(8,8).Value = 1792.124268
$ws.Cells.Item(8,9).Value = 0.8773093221949784
$ws.Cells.Item(8,10).Value = 0.8773093221949785
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.2909853333333334
$ws.Cells.Item(8,14).Value = 0.8729560000000001
$ws.Cells.Item(8,15).Value = 0.03275300561492853
$ws.Cells.Item(8,16).Value = 0.03275300561492853
$ws.Cells.Item(8,17).Value = 173.8272924995787
$ws.Cells.Item(8,18).Value = 1564.445632496208
$ws.Cells.Item(8,19).Value = 0.02873451715588127
$ws.Cells.Item(8,20).Value = 0.02873451715588127
# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Nid1"
$ws.Cells.Item(9,3).Value = "Itgb3"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 597.374756
$ws.Cells.Item(9,8).Value = 1792.124268
$ws.Cells.Item(9,9).Value = 0.8773093221949784
$ws.Cells.Item(9,10).Value = 0.8773093221949785
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.4814053333333333
$ws.Cells.Item(9,14).Value = 1.444216
$ws.Cells.Item(9,15).Value = 0.0541864822020464
$ws.Cells.Item(9,16).Value = 0.05418648220204641
$ws.Cells.Item(9,17).Value = 287.5793935370987
$ws.Cells.Item(9,18).Value = 2588.214541833888
$ws.Cells.Item(9,19).Value = 0.04753830597280759
$ws.Cells.Item(9,20).Value = 0.0475383059728076
# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Nid1"
$ws.Cells.Item(10,3).Value = "Itgb3"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.1960536666666667
$ws.Cells.Item(10,8).Value = 0.5881609999999999
$ws.Cells.Item(10,9).Value = 0.0002879259755950811
$ws.Cells.Item(10,10).Value = 0.0002879259755950811
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.778439
$ws.Cells.Item(10,14).Value = 11.335317
$ws.Cells.Item(10,15).Value = 0.4252971528324392
$ws.Cells.Item(10,16).Value = 0.4252971528324392
$ws.Cells.Item(10,17).Value = 0.7407768202263333
$ws.Cells.Item(10,18).Value = 6.666991382036999
$ws.Cells.Item(10,19).Value = 0.0001224540976470904
$ws.Cells.Item(10,20).Value = 0.0001224540976470904
# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Nid1"
$ws.Cells.Item(11,3).Value = "Itgb3"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.1960536666666667
$ws.Cells.Item(11,8).Value = 0.5881609999999999
$ws.Cells.Item(11,9).Value = 0.0002879259755950811
$ws.Cells.Item(11,10).Value = 0.0002879259755950811
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 4.333403333333333
$ws.Cells.Item(11,14).Value = 13.00021
$ws.Cells.Item(11,15).Value = 0.4877633593505858
$ws.Cells.Item(11,16).Value = 0.4877633593505858
$ws.Cells.Item(11,17).Value = 0.8495796126455555
$ws.Cells.Item(11,18).Value = 7.646216513809999
$ws.Cells.Item(11,19).Value = 0.0001404397411005515
$ws.Cells.Item(11,20).Value = 0.0001404397411005516
# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Nid1"
$ws.Cells.Item(12,3).Value = "Itgb3"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 0.6666666666666666
$ws.Cells.Item(12,7).Value = 0.1960536666666667
$ws.Cells.Item(12,8).Value = 0.5881609999999999
$ws.Cells.Item(12,9).Value = 0.0002879259755950811
$ws.Cells.Item(12,10).Value = 0.0002879259755950811
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.2909853333333334
$ws.Cells.Item(12,14).Value = 0.8729560000000001
$ws.Cells.Item(12,15).Value = 0.03275300561492853
$ws.Cells.Item(12,16).Value = 0.03275300561492853
$ws.Cells.Item(12,17).Value = 0.05704874154622223
$ws.Cells.Item(12,18).Value = 0.513438673916
$ws.Cells.Item(12,19).Value = 0.000009430441095349466
$ws.Cells.Item(12,20).Value = 0.000009430441095349468
# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Nid1"
$ws.Cells.Item(13,3).Value = "Itgb3"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 0.6666666666666666
$ws.Cells.Item(13,7).Value = 0.1960536666666667
$ws.Cells.Item(13,8).Value = 0.5881609999999999
$ws.Cells.Item(13,9).Value = 0.0002879259755950811
$ws.Cells.Item(13,10).Value = 0.0002879259755950811
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.4814053333333333
$ws.Cells.Item(13,14).Value = 1.444216
$ws.Cells.Item(13,15).Value = 0.0541864822020464
$ws.Cells.Item(13,16).Value = 0.05418648220204641
$ws.Cells.Item(13,17).Value = 0.09438128075288887
$ws.Cells.Item(13,18).Value = 0.8494315267759999
$ws.Cells.Item(13,19).Value = 0.00001560169575208971
$ws.Cells.Item(13,20).Value = 0.00001560169575208971
# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Nid1"
$ws.Cells.Item(14,3).Value = "Itgb3"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 27.92897333333333
$ws.Cells.Item(14,8).Value = 83.78692
$ws.Cells.Item(14,9).Value = 0.04101671257207978
$ws.Cells.Item(14,10).Value = 0.04101671257207978
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 3.778439
$ws.Cells.Item(14,14).Value = 11.335317
$ws.Cells.Item(14,15).Value = 0.4252971528324392
$ws.Cells.Item(14,16).Value = 0.4252971528324392
$ws.Cells.Item(14,17).Value = 105.5279220726267
$ws.Cells.Item(14,18).Value = 949.7512986536399
$ws.Cells.Item(14,19).Value = 0.01744429107545204
$ws.Cells.Item(14,20).Value = 0.01744429107545205
# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Nid1"
$ws.Cells.Item(15,3).Value = "Itgb3"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 27.92897333333333
$ws.Cells.Item(15,8).Value = 83.78692
$ws.Cells.Item(15,9).Value = 0.04101671257207978
$ws.Cells.Item(15,10).Value = 0.04101671257207978
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 4.333403333333333
$ws.Cells.Item(15,14).Value = 13.00021
$ws.Cells.Item(15,15).Value = 0.4877633593505858
$ws.Cells.Item(15,16).Value = 0.4877633593505858
$ws.Cells.Item(15,17).Value = 121.0275061392444
$ws.Cells.Item(15,18).Value = 1089.2475552532
$ws.Cells.Item(15,19).Value = 0.02000644951367504
$ws.Cells.Item(15,20).Value = 0.02000644951367504
# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Nid1"
$ws.Cells.Item(16,3).Value = "Itgb3"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 27.92897333333333
$ws.Cells.Item(16,8).Value = 83.78692
$ws.Cells.Item(16,9).Value = 0.04101671257207978
$ws.Cells.Item(16,10).Value = 0.04101671257207978
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.2909853333333334
$ws.Cells.Item(16,14).Value = 0.8729560000000001
$ws.Cells.Item(16,15).Value = 0.03275300561492853
$ws.Cells.Item(16,16).Value = 0.03275300561492853
$ws.Cells.Item(16,17).Value = 8.126921615057778
$ws.Cells.Item(16,18).Value = 73.14229453552
$ws.Cells.Item(16,19).Value = 0.001343420617179239
$ws.Cells.Item(16,20).Value = 0.001343420617179239
# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Nid1"
$ws.Cells.Item(17,3).Value = "Itgb3"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 27.92897333333333
$ws.Cells.Item(17,8).Value = 83.78692
$ws.Cells.Item(17,9).Value = 0.04101671257207978
$ws.Cells.Item(17,10).Value = 0.04101671257207978
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 0.4814053333333333
$ws.Cells.Item(17,14).Value = 1.444216
$ws.Cells.Item(17,15).Value = 0.0541864822020464
$ws.Cells.Item(17,16).Value = 0.05418648220204641
$ws.Cells.Item(17,17).Value = 13.44515671719111
$ws.Cells.Item(17,18).Value = 121.00641045472
$ws.Cells.Item(17,19).Value = 0.002222551365773454
$ws.Cells.Item(17,20).Value = 0.002222551365773454
